$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC - Cadastro Visitante")
$ws.Activate()
$win = $excel.ActiveWindow
$pane = $win.ActivePane
try {
  $pane.ScrollRow = 15
  $pane.ScrollColumn = 1
  Write-Host "Pane scroll set ok; ScrollRow=" $pane.ScrollRow
} catch {
  Write-Host "Pane scroll failed: $_"
}
$ws.Range("A25").Select()
